$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Content.Find.Execute("2025-10-11 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-12 Sunday", 2) | Out-Null

# Update the 20x5 table of arithmetic answers, cell by cell (row, col are 1-indexed)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92-78=14"
$t.Cell(1, 2).Range.Text = "17+17=34"
$t.Cell(1, 3).Range.Text = "75-16=59"
$t.Cell(1, 4).Range.Text = "32-15=17"
$t.Cell(1, 5).Range.Text = "22-19=3"

$t.Cell(2, 1).Range.Text = "36+48=84"
$t.Cell(2, 2).Range.Text = "62-28=34"
$t.Cell(2, 3).Range.Text = "34+8=42"
$t.Cell(2, 4).Range.Text = "49+24=73"
$t.Cell(2, 5).Range.Text = "78+3=81"

$t.Cell(3, 1).Range.Text = "33-15=18"
$t.Cell(3, 2).Range.Text = "66+29=95"
$t.Cell(3, 3).Range.Text = "91-6=85"
$t.Cell(3, 4).Range.Text = "77+15=92"
$t.Cell(3, 5).Range.Text = "35+38=73"

$t.Cell(4, 1).Range.Text = "7+48=55"
$t.Cell(4, 2).Range.Text = "24+59=83"
$t.Cell(4, 3).Range.Text = "53+28=81"
$t.Cell(4, 4).Range.Text = "82-57=25"
$t.Cell(4, 5).Range.Text = "28+47=75"

$t.Cell(5, 1).Range.Text = "86-77=9"
$t.Cell(5, 2).Range.Text = "18+28=46"
$t.Cell(5, 3).Range.Text = "91-54=37"
$t.Cell(5, 4).Range.Text = "88+8=96"
$t.Cell(5, 5).Range.Text = "44-6=38"

$t.Cell(6, 1).Range.Text = "48+45=93"
$t.Cell(6, 2).Range.Text = "58+15=73"
$t.Cell(6, 3).Range.Text = "47+44=91"
$t.Cell(6, 4).Range.Text = "50-36=14"
$t.Cell(6, 5).Range.Text = "89+5=94"

$t.Cell(7, 1).Range.Text = "71-48=23"
$t.Cell(7, 2).Range.Text = "60-4=56"
$t.Cell(7, 3).Range.Text = "8+64=72"
$t.Cell(7, 4).Range.Text = "27+66=93"
$t.Cell(7, 5).Range.Text = "17+27=44"

$t.Cell(8, 1).Range.Text = "26+36=62"
$t.Cell(8, 2).Range.Text = "81-6=75"
$t.Cell(8, 3).Range.Text = "24+18=42"
$t.Cell(8, 4).Range.Text = "39+48=87"
$t.Cell(8, 5).Range.Text = "47+24=71"

$t.Cell(9, 1).Range.Text = "41-19=22"
$t.Cell(9, 2).Range.Text = "75-16=59"
$t.Cell(9, 3).Range.Text = "77+15=92"
$t.Cell(9, 4).Range.Text = "56+7=63"
$t.Cell(9, 5).Range.Text = "18+67=85"

$t.Cell(10, 1).Range.Text = "71-13=58"
$t.Cell(10, 2).Range.Text = "80-57=23"
$t.Cell(10, 3).Range.Text = "86-37=49"
$t.Cell(10, 4).Range.Text = "25+57=82"
$t.Cell(10, 5).Range.Text = "80-14=66"

$t.Cell(11, 1).Range.Text = "90-39=51"
$t.Cell(11, 2).Range.Text = "54-25=29"
$t.Cell(11, 3).Range.Text = "65-46=19"
$t.Cell(11, 4).Range.Text = "76-19=57"
$t.Cell(11, 5).Range.Text = "62-35=27"

$t.Cell(12, 1).Range.Text = "29+37=66"
$t.Cell(12, 2).Range.Text = "52-45=7"
$t.Cell(12, 3).Range.Text = "61-38=23"
$t.Cell(12, 4).Range.Text = "66+28=94"
$t.Cell(12, 5).Range.Text = "89+3=92"

$t.Cell(13, 1).Range.Text = "21-3=18"
$t.Cell(13, 2).Range.Text = "95-86=9"
$t.Cell(13, 3).Range.Text = "22-5=17"
$t.Cell(13, 4).Range.Text = "42-28=14"
$t.Cell(13, 5).Range.Text = "75-27=48"

$t.Cell(14, 1).Range.Text = "60-59=1"
$t.Cell(14, 2).Range.Text = "94-25=69"
$t.Cell(14, 3).Range.Text = "6+88=94"
$t.Cell(14, 4).Range.Text = "28+63=91"
$t.Cell(14, 5).Range.Text = "61-59=2"

$t.Cell(15, 1).Range.Text = "8+29=37"
$t.Cell(15, 2).Range.Text = "9+24=33"
$t.Cell(15, 3).Range.Text = "41-24=17"
$t.Cell(15, 4).Range.Text = "80-15=65"
$t.Cell(15, 5).Range.Text = "40-28=12"

$t.Cell(16, 1).Range.Text = "15+59=74"
$t.Cell(16, 2).Range.Text = "25+39=64"
$t.Cell(16, 3).Range.Text = "4+7=11"
$t.Cell(16, 4).Range.Text = "85-36=49"
$t.Cell(16, 5).Range.Text = "76+6=82"

$t.Cell(17, 1).Range.Text = "62-8=54"
$t.Cell(17, 2).Range.Text = "44+7=51"
$t.Cell(17, 3).Range.Text = "30-14=16"
$t.Cell(17, 4).Range.Text = "39+3=42"
$t.Cell(17, 5).Range.Text = "46+5=51"

$t.Cell(18, 1).Range.Text = "92-23=69"
$t.Cell(18, 2).Range.Text = "54+39=93"
$t.Cell(18, 3).Range.Text = "48+17=65"
$t.Cell(18, 4).Range.Text = "76+8=84"
$t.Cell(18, 5).Range.Text = "36+59=95"

$t.Cell(19, 1).Range.Text = "13+59=72"
$t.Cell(19, 2).Range.Text = "92-55=37"
$t.Cell(19, 3).Range.Text = "38+34=72"
$t.Cell(19, 4).Range.Text = "57+29=86"
$t.Cell(19, 5).Range.Text = "23-17=6"

$t.Cell(20, 1).Range.Text = "43-34=9"
$t.Cell(20, 2).Range.Text = "62-46=16"
$t.Cell(20, 3).Range.Text = "66+17=83"
$t.Cell(20, 4).Range.Text = "78+18=96"
$t.Cell(20, 5).Range.Text = "4+18=22"
